$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two workers who left (rows shift up automatically) ---
# Row 30 = "Rodd McGowan", Row 24 = "Mason Gardner" (delete bottom-most first
# so the earlier row index used for the second delete stays valid).
$ws.Rows("30").Delete()
$ws.Rows("24").Delete()

# --- Fill in newly collected scanner ID / S-N data for a few workers ---
# (row numbers below are the POST-deletion row numbers)
$ws.Range("C27").Value = "7A21E43A"
$ws.Range("D27").Value = "N521D5060004"
$ws.Range("C27").HorizontalAlignment = -4131   # xlLeft - match the row's existing look

$ws.Range("C16").Value = "0E9D9A52"
$ws.Range("D16").Value = "N521D5060011"

$ws.Range("C17").Value = "BC590752"
$ws.Range("D17").Value = "N521D5060007"

$ws.Range("C34").Value = "FECB7890"
$ws.Range("D34").Value = "N521D5060010"

# --- Turn the data range into a proper Excel Table named "WorkerList" ---
$ws.AutoFilterMode = $false
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:D35"), $null, 1)
$tbl.Name = "WorkerList"

# --- Tidy up the defined name that used to back the manual AutoFilter ---
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$4:`$A`$52"

# --- Refresh the stale sort-state metadata to match the new data extent ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add2($ws.Range("A1:A48"))
$ws.Sort.SetRange($ws.Range("A4:C48"))
$ws.Sort.Apply()

# --- Minor formatting / view tweaks ---
$ws.Columns("B").ColumnWidth = 9.6
$ws.Range("C5").Select()
